$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 0
    3  = 0
    4  = 0
    6  = 0
    7  = 0
    9  = 11186.53241152766
    10 = 11186.53241152765
    11 = 162.6848836413069
    12 = 9234.31380783197
    13 = 162.6848836413069
    17 = 130147.9069130455
    19 = 143934.1102782082
    20 = 143934.1102782082
    21 = 1952.218603695683
    24 = 30801.87040944644
    25 = 30801.87040944641
    28 = 616.0374081889174
    29 = 616.0374081889174
    30 = 30801.87040944641
    35 = 20217.59999999929
    36 = 20217.59999999929
    38 = -6138.278472965609
    39 = -6138.278472965603
    41 = 6138.278472965603
    42 = 122765.5694593122
    43 = 122765.5694593122
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}
